$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '41.092.25'
$ws.Range("E2").Value = '  +0.87%  '

Set-TextValue "D3" '2.177.55'
$ws.Range("E3").Value = '  -0.46%  '

$ws.Range("E4").Value = '  +0.14%  '

Set-TextValue "D5" '254.47'
$ws.Range("E5").Value = '  +5.57%  '

Set-TextValue "D6" '0.627'
$ws.Range("E6").Value = '  +1.14%  '

Set-TextValue "D7" '67.79'
$ws.Range("E7").Value = '  -1.25%  '

$ws.Range("E8").Value = '  +0.01%  '

Set-TextValue "D9" '0.574'
$ws.Range("E9").Value = '  +7.20%  '

Set-TextValue "D10" '37.58'
$ws.Range("E10").Value = '  +3.23%  '

$ws.Range("E11").Value = '  +2.39%  '

Set-TextValue "D12" '0.0930'
$ws.Range("E12").Value = '  -1.39%  '

Set-TextValue "D13" '7.10'
$ws.Range("E13").Value = '  +8.34%  '

Set-TextValue "D14" '0.104'
$ws.Range("E14").Value = '  +0.68%  '

Set-TextValue "D15" '2.504.55'
$ws.Range("E15").Value = '  -0.34%  '

Set-TextValue "D16" '0.870'
$ws.Range("E16").Value = '  +5.15%  '

Set-TextValue "D17" '14.44'
$ws.Range("E17").Value = '  -1.08%  '

Set-TextValue "D18" '2.180.21'
$ws.Range("E18").Value = '  -0.34%  '

Set-TextValue "D19" '41.165.52'
$ws.Range("E19").Value = '  +1.23%  '

Set-TextValue "D20" '0.0₃0952'
$ws.Range("E20").Value = '  +1.64%  '

Set-TextValue "D21" '6.15'
$ws.Range("E21").Value = '  +2.07%  '

Set-TextValue "D22" '71.63'
$ws.Range("E22").Value = '  -0.89%  '

Set-TextValue "D23" '231.63'
$ws.Range("E23").Value = '  +1.06%  '

Set-TextValue "D24" '2.03'
$ws.Range("E24").Value = '  +0.75%  '

Set-TextValue "D25" '3.94'
$ws.Range("E25").Value = '  +10.19%  '

Set-TextValue "D26" '11.80'
$ws.Range("E26").Value = '  +22.38%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D27" '1.00'
$ws.Range("E27").Value = '  +0.04%  '

$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D28" '2.53'
$ws.Range("E28").Value = '  +5.92%  '

Set-TextValue "D29" '2.17'
$ws.Range("E29").Value = '  +0.21%  '

Set-TextValue "D30" '168.25'
$ws.Range("E30").Value = '  +0.01%  '

Set-TextValue "D31" '20.59'
$ws.Range("E31").Value = '  +2.12%  '

Set-TextValue "D32" '0.116'
$ws.Range("E32").Value = '  -0.81%  '

Set-TextValue "D33" '0.0746'
$ws.Range("E33").Value = '  +7.11%  '

Set-TextValue "D34" '0.123'
$ws.Range("E34").Value = '  +0.52%  '

Set-TextValue "D35" '5.42'
$ws.Range("E35").Value = '  +6.37%  '

Set-TextValue "D36" '26.71'
$ws.Range("E36").Value = '  +14.83%  '

$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D37" '4.62'
$ws.Range("E37").Value = '  +1.48%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D38" '4.11'
$ws.Range("E38").Value = '  +7.69%  '

$ws.Range("E39").Value = '  +13.25%  '

Set-TextValue "D40" '2.20'
$ws.Range("E40").Value = '  -2.47%  '

Set-TextValue "D41" '12.52'
$ws.Range("E41").Value = '  +22.26%  '

Set-TextValue "D42" '5.68'
$ws.Range("E42").Value = '  -1.60%  '

Set-TextValue "D43" '64.28'
$ws.Range("E43").Value = '  +2.42%  '

Set-TextValue "D44" '5.06'
$ws.Range("E44").Value = '  +5.04%  '

Set-TextValue "D45" '0.200'
$ws.Range("E45").Value = '  +5.24%  '

Set-TextValue "D46" '8.64'
$ws.Range("E46").Value = '  +1.06%  '

Set-TextValue "D47" '0.100'
$ws.Range("E47").Value = '  +2.79%  '

$ws.Range("E48").Value = '  +0.51%  '

Set-TextValue "D49" '1.13'
$ws.Range("E49").Value = '  +4.54%  '

Set-TextValue "D50" '1.17'
$ws.Range("E50").Value = '  +1.32%  '

Set-TextValue "D51" '4.28'
$ws.Range("E51").Value = '  -4.55%  '
